# netCrypto.xlsx - "Add files via upload"
#
# The commit re-saves the daily-updates workbook after the user scrolled to
# a new row and refreshed the USD Amount total:
#   - the workbook window was nudged back towards the left edge of the screen
#   - the active selection on SheetName1 moved from T3 to T2
#   - T2 (USD Amount for the Deposit/Crypto/Roobic row) was updated
#     from 73959 to 105329

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook window was moved horizontally (xWindow -28920 -> -120); the
# vertical position (yWindow) and size were left untouched.
$excel.ActiveWindow.Left = -120

# Selection moved from T3 to T2.
$ws.Range("T2").Select()

# USD Amount total refreshed: 73959 -> 105329.
$ws.Range("T2").Value = 105329
